# "add OILK to sectors"
#
# The sheet lists tickers (with min/max weight constraints) grouped in
# blocks (bonds, broad equity, sectors, commodities). The sector block
# ends with the commodities entry "GLD" on row 19. We need to insert a
# new sector row for "OILK" (same 0 / 0.01 min/max weight constraints as
# GLD) directly above it, pushing GLD down to row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Investments")

# Copy row 18 ("EEM") - it carries the same formatting/styles that the
# new row should have (bold/left-bordered ticker cell) - then insert a
# new, fully-formatted row above row 19 ("GLD"), shifting GLD (and
# nothing else, since it's the last data row) down to row 20.
$ws.Rows.Item(18).Copy() | Out-Null
$ws.Rows.Item(19).Insert() | Out-Null
$excel.CutCopyMode = $false

# Fill in the new row with the OILK data.
$ws.Range("A19").Value = "OILK"
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0.01
